# Update cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "42.964.43"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.330.83"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'303.15"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'96.00"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D10").Value = "'34.18"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("D11").Value = "'19.20"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").Value = "'0.0786"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "2.693.38"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "2.340.01"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").Value = "42.920.13"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'236.78"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D27").Value = "'24.63"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  -13.69%  "
$ws.Range("D29").Value = "'9.15"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").Value = "'31.61"
$ws.Range("E30").Value = "  -3.71%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'138.54"
$ws.Range("E32").Value = "  -16.51%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "'17.75"
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("D35").Value = "'0.0704"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").Value = "'4.41"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.75"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'22.20"
$ws.Range("E41").Value = "  +22.71%  "
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "1.933.02"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "'10.03"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").Value = "2.561.22"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "'53.67"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'73.02"
$ws.Range("E51").Value = "  +1.74%  "
